# Auto-generated edit script: update market price / profit columns (H:N)
# across multiple sheets, reflecting a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1369.8889
$ws.Range("I28").Value = 543.625
$ws.Range("J28").Value = 7980
$ws.Range("K28").Value = 543.625
$ws.Range("L28").Value = 7980
$ws.Range("M28").Value = -58.625
$ws.Range("N28").Value = -8950

$ws.Range("H43").Value = 7311.9414
$ws.Range("I43").Value = 3250
$ws.Range("J43").Value = 7853.533
$ws.Range("K43").Value = 3250
$ws.Range("L43").Value = 7853.533
$ws.Range("M43").Value = -3181
$ws.Range("N43").Value = -7991.533

$ws.Range("H88").Value = 4658.5557
$ws.Range("I88").Value = 2716.6667
$ws.Range("J88").Value = 5629.5
$ws.Range("K88").Value = 2716.6667
$ws.Range("L88").Value = 5629.5
$ws.Range("M88").Value = -2310.6667
$ws.Range("N88").Value = -6441.5

$ws.Range("H91").Value = 4658.5557
$ws.Range("I91").Value = 2716.6667
$ws.Range("J91").Value = 5629.5
$ws.Range("K91").Value = 2716.6667
$ws.Range("L91").Value = 5629.5
$ws.Range("M91").Value = -1312.6667
$ws.Range("N91").Value = -8437.5

$ws.Range("H98").Value = 3985
$ws.Range("I98").Value = 1279
$ws.Range("J98").Value = 10750
$ws.Range("K98").Value = 1279
$ws.Range("L98").Value = 10750
$ws.Range("M98").Value = 219
$ws.Range("N98").Value = -13746

$ws.Range("H122").Value = 3985
$ws.Range("I122").Value = 1279
$ws.Range("J122").Value = 10750
$ws.Range("K122").Value = 3837
$ws.Range("L122").Value = 32250
$ws.Range("M122").Value = -1387
$ws.Range("N122").Value = -37150

$ws.Range("H132").Value = 3289.3447
$ws.Range("I132").Value = 1680.5532
$ws.Range("J132").Value = 10163.272
$ws.Range("K132").Value = 5041.6596
$ws.Range("L132").Value = 30489.816
$ws.Range("M132").Value = -2511.6596
$ws.Range("N132").Value = -35549.81600000001

$ws.Range("H137").Value = 4247.864
$ws.Range("I137").Value = 4452.9443
$ws.Range("J137").Value = 3325
$ws.Range("K137").Value = 13358.8329
$ws.Range("L137").Value = 9975
$ws.Range("M137").Value = -10808.8329
$ws.Range("N137").Value = -15075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1683.2
$ws.Range("I22").Value = 805.3333
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 805.3333
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -506.3333
$ws.Range("N22").Value = -3598

$ws.Range("H28").Value = 40259.5
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 40259.5
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 40259.5
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -40643.5

$ws.Range("H39").Value = 900
$ws.Range("I39").Value = 900
$ws.Range("K39").Value = 900
$ws.Range("M39").Value = -380

$ws.Range("H41").Value = 2556
$ws.Range("I41").Value = 2556
$ws.Range("K41").Value = 2556
$ws.Range("M41").Value = -2142

$ws.Range("I61").Value = 356612.16
$ws.Range("J61").Value = 479788.47
$ws.Range("K61").Value = 356612.16
$ws.Range("L61").Value = 479788.47
$ws.Range("M61").Value = -356400.16
$ws.Range("N61").Value = -480212.47

$ws.Range("H99").Value = 40259.5
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 40259.5
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 40259.5
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -46249.5

$ws.Range("H132").Value = 20888.447
$ws.Range("I132").Value = 27023.684
$ws.Range("K132").Value = 81071.052
$ws.Range("M132").Value = -78541.052

$ws.Range("I136").Value = 356612.16
$ws.Range("J136").Value = 479788.47
$ws.Range("K136").Value = 1069836.48
$ws.Range("L136").Value = 1439365.41
$ws.Range("M136").Value = -1067286.48
$ws.Range("N136").Value = -1444465.41

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 40
$ws.Range("N8").ClearContents()

$ws.Range("H20").Value = 1240.5714
$ws.Range("I20").Value = 1144.409
$ws.Range("J20").Value = 1593.1666
$ws.Range("K20").Value = 1144.409
$ws.Range("L20").Value = 1593.1666
$ws.Range("M20").Value = -897.4090000000001
$ws.Range("N20").Value = -2087.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6464.8184
$ws.Range("I58").Value = 11102.6
$ws.Range("K58").Value = 11102.6
$ws.Range("M58").Value = -10899.6

$ws.Range("H99").Value = 49483.523
$ws.Range("I99").Value = 63925.625
$ws.Range("J99").Value = 3268.8
$ws.Range("K99").Value = 63925.625
$ws.Range("L99").Value = 3268.8
$ws.Range("M99").Value = -62427.625
$ws.Range("N99").Value = -6264.8

$ws.Range("H126").Value = 49483.523
$ws.Range("I126").Value = 63925.625
$ws.Range("J126").Value = 3268.8
$ws.Range("K126").Value = 191776.875
$ws.Range("L126").Value = 9806.400000000001
$ws.Range("M126").Value = -189306.875
$ws.Range("N126").Value = -14746.4

$ws.Range("H136").Value = 6464.8184
$ws.Range("I136").Value = 11102.6
$ws.Range("K136").Value = 33307.8
$ws.Range("M136").Value = -30757.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 894.5
$ws.Range("I5").Value = 894.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2683.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2571.5
$ws.Range("N5").ClearContents()

$ws.Range("H17").Value = 846.4286
$ws.Range("I17").Value = 111
$ws.Range("J17").Value = 1398
$ws.Range("K17").Value = 333
$ws.Range("L17").Value = 4194
$ws.Range("M17").Value = -164
$ws.Range("N17").Value = -4532

$ws.Range("H39").Value = 4080.8
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4080.8
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 12242.4
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -12830.4

$ws.Range("H132").Value = 8638.666999999999
$ws.Range("I132").Value = 1734.5
$ws.Range("J132").Value = 14162
$ws.Range("K132").Value = 15610.5
$ws.Range("L132").Value = 127458
$ws.Range("M132").Value = -13080.5
$ws.Range("N132").Value = -132518

$ws.Range("H135").Value = 894.5
$ws.Range("I135").Value = 894.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8050.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -5515.5
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 43149.5
$ws.Range("J128").Value = 43149.5
$ws.Range("L128").Value = 43149.5
$ws.Range("N128").Value = -53109.5

$ws.Range("H132").Value = 5791.304
$ws.Range("I132").Value = 9181.875
$ws.Range("J132").Value = 3983
$ws.Range("K132").Value = 27545.625
$ws.Range("L132").Value = 11949
$ws.Range("M132").Value = -25015.625
$ws.Range("N132").Value = -17009

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2213.1738
$ws.Range("I122").Value = 2080.2
$ws.Range("K122").Value = 6240.599999999999
$ws.Range("M122").Value = -3790.599999999999

$ws.Range("H136").Value = 3806.1965
$ws.Range("I136").Value = 2376.389
$ws.Range("K136").Value = 7129.167
$ws.Range("M136").Value = -4579.167

$ws.Range("H139").Value = 42500
$ws.Range("J139").Value = 42500
$ws.Range("L139").Value = 42500
$ws.Range("N139").Value = -52780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2337.0286
$ws.Range("I132").Value = 1382.8096
$ws.Range("J132").Value = 3768.3572
$ws.Range("K132").Value = 4148.4288
$ws.Range("L132").Value = 11305.0716
$ws.Range("M132").Value = -1618.4288
$ws.Range("N132").Value = -16365.0716
